# away6+hours.xlsx - "data up to 17th"
#   - backfills the previously-missing American Samoa (col E) values on
#     rows 95-101
#   - backfills the rest of the state columns (B:D, F:BE) for the rows
#     dated 06-10 Aug 2020 that already existed with only a date label
#     in col A (rows 189-193)
#   - adds two brand-new fully populated rows for 11 and 12 Aug 2020
#     (rows 194-195)
#   - adds five new date-only placeholder rows for 13-17 Aug 2020
#     (rows 196-200), matching the same pattern the workbook already
#     used for not-yet-reported dates

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 95
$ws.Cells.Item(95, 5).Value = 0.10633903133903

# Row 96
$ws.Cells.Item(96, 5).Value = 0.1143475572047

# Row 97
$ws.Cells.Item(97, 5).Value = 0.12175324675325

# Row 98
$ws.Cells.Item(98, 5).Value = 0.14635225885226

# Row 99
$ws.Cells.Item(99, 5).Value = 0.13388888888889

# Row 100
$ws.Cells.Item(100, 5).Value = 0.16659673659674

# Row 101
$ws.Cells.Item(101, 5).Value = 0.1247923312068

# Row 189
$ws.Cells.Item(189, 2).Value = 0.050599845851136
$ws.Cells.Item(189, 3).Value = 0.05737187068404
$ws.Cells.Item(189, 4).Value = 0.047005831429903
$ws.Cells.Item(189, 6).Value = 0.043146103170363
$ws.Cells.Item(189, 7).Value = 0.046638236245796
$ws.Cells.Item(189, 8).Value = 0.050441598304295
$ws.Cells.Item(189, 9).Value = 0.0439840881723
$ws.Cells.Item(189, 10).Value = 0.046806704259656
$ws.Cells.Item(189, 11).Value = 0.044448796698912
$ws.Cells.Item(189, 12).Value = 0.046528870188271
$ws.Cells.Item(189, 13).Value = 0.050307370472125
$ws.Cells.Item(189, 14).Value = 0.057413097782474
$ws.Cells.Item(189, 15).Value = 0.031001433419907
$ws.Cells.Item(189, 16).Value = 0.045654679208709
$ws.Cells.Item(189, 17).Value = 0.049821844993392
$ws.Cells.Item(189, 18).Value = 0.043247276004654
$ws.Cells.Item(189, 19).Value = 0.049973610373662
$ws.Cells.Item(189, 20).Value = 0.044399428131163
$ws.Cells.Item(189, 21).Value = 0.045885839977657
$ws.Cells.Item(189, 22).Value = 0.054360510211892
$ws.Cells.Item(189, 23).Value = 0.046004787263399
$ws.Cells.Item(189, 24).Value = 0.041591488030797
$ws.Cells.Item(189, 25).Value = 0.04721995078091
$ws.Cells.Item(189, 26).Value = 0.049275315020692
$ws.Cells.Item(189, 27).Value = 0.046905493283737
$ws.Cells.Item(189, 28).Value = 0.052678432214427
$ws.Cells.Item(189, 29).Value = 0.064677237308593
$ws.Cells.Item(189, 30).Value = 0.06015423689422
$ws.Cells.Item(189, 31).Value = 0.060469309645183
$ws.Cells.Item(189, 32).Value = 0.047967647115805
$ws.Cells.Item(189, 33).Value = 0.044339046838329
$ws.Cells.Item(189, 34).Value = 0.053294399635864
$ws.Cells.Item(189, 35).Value = 0.052297494433935
$ws.Cells.Item(189, 36).Value = 0.048822788416176
$ws.Cells.Item(189, 37).Value = 0.05049552715532
$ws.Cells.Item(189, 38).Value = 0.050432207723194
$ws.Cells.Item(189, 39).Value = 0.048888445792371
$ws.Cells.Item(189, 40).Value = 0.048645943823381
$ws.Cells.Item(189, 41).Value = 0.052239848033756
$ws.Cells.Item(189, 42).Value = 0.048317080967942
$ws.Cells.Item(189, 43).Value = 0.050096067448388
$ws.Cells.Item(189, 44).Value = 0.083819966287892
$ws.Cells.Item(189, 45).Value = 0.055306609688787
$ws.Cells.Item(189, 46).Value = 0.053826754440282
$ws.Cells.Item(189, 47).Value = 0.052793853771579
$ws.Cells.Item(189, 48).Value = 0.053502724498918
$ws.Cells.Item(189, 49).Value = 0.047154402639104
$ws.Cells.Item(189, 50).Value = 0.057894781009057
$ws.Cells.Item(189, 51).Value = 0.049463190349083
$ws.Cells.Item(189, 52).Value = 0.08571278696910101
$ws.Cells.Item(189, 53).Value = 0.053146487325292
$ws.Cells.Item(189, 54).Value = 0.048457504661775
$ws.Cells.Item(189, 55).Value = 0.054981403294058
$ws.Cells.Item(189, 56).Value = 0.057290174120725
$ws.Cells.Item(189, 57).Value = 0.052538238067564

# Row 190
$ws.Cells.Item(190, 2).Value = 0.051252978723225
$ws.Cells.Item(190, 3).Value = 0.043514127562904
$ws.Cells.Item(190, 4).Value = 0.041363011566801
$ws.Cells.Item(190, 6).Value = 0.043603069701841
$ws.Cells.Item(190, 7).Value = 0.044928184634195
$ws.Cells.Item(190, 8).Value = 0.051647905848707
$ws.Cells.Item(190, 9).Value = 0.047323929877081
$ws.Cells.Item(190, 10).Value = 0.047319119856335
$ws.Cells.Item(190, 11).Value = 0.044369658793154
$ws.Cells.Item(190, 12).Value = 0.043776130730086
$ws.Cells.Item(190, 13).Value = 0.049155927133191
$ws.Cells.Item(190, 14).Value = 0.053158911477339
$ws.Cells.Item(190, 15).Value = 0.034679551947387
$ws.Cells.Item(190, 16).Value = 0.040658638486645
$ws.Cells.Item(190, 17).Value = 0.050165785822023
$ws.Cells.Item(190, 18).Value = 0.046377563655325
$ws.Cells.Item(190, 19).Value = 0.04822711986515
$ws.Cells.Item(190, 20).Value = 0.038850432876626
$ws.Cells.Item(190, 21).Value = 0.039108655597086
$ws.Cells.Item(190, 22).Value = 0.045790768591373
$ws.Cells.Item(190, 23).Value = 0.037662529855784
$ws.Cells.Item(190, 24).Value = 0.036016208628749
$ws.Cells.Item(190, 25).Value = 0.039012113841025
$ws.Cells.Item(190, 26).Value = 0.037372126448531
$ws.Cells.Item(190, 27).Value = 0.032428918528143
$ws.Cells.Item(190, 28).Value = 0.037171535128803
$ws.Cells.Item(190, 29).Value = 0.072286331938126
$ws.Cells.Item(190, 30).Value = 0.044493742378422
$ws.Cells.Item(190, 31).Value = 0.040952479617276
$ws.Cells.Item(190, 32).Value = 0.037916354471223
$ws.Cells.Item(190, 33).Value = 0.03565859461879
$ws.Cells.Item(190, 34).Value = 0.037391642300712
$ws.Cells.Item(190, 35).Value = 0.034678773157402
$ws.Cells.Item(190, 36).Value = 0.033317641656494
$ws.Cells.Item(190, 37).Value = 0.038150267793167
$ws.Cells.Item(190, 38).Value = 0.036437180526116
$ws.Cells.Item(190, 39).Value = 0.037897875803645
$ws.Cells.Item(190, 40).Value = 0.041031797670719
$ws.Cells.Item(190, 41).Value = 0.042060133319046
$ws.Cells.Item(190, 42).Value = 0.039597309156181
$ws.Cells.Item(190, 43).Value = 0.04182385142185
$ws.Cells.Item(190, 44).Value = 0.071933318481847
$ws.Cells.Item(190, 45).Value = 0.042147086942694
$ws.Cells.Item(190, 46).Value = 0.04068620297373
$ws.Cells.Item(190, 47).Value = 0.040837990534141
$ws.Cells.Item(190, 48).Value = 0.046223773205754
$ws.Cells.Item(190, 49).Value = 0.038859849848363
$ws.Cells.Item(190, 50).Value = 0.047748783338324
$ws.Cells.Item(190, 51).Value = 0.039906295350177
$ws.Cells.Item(190, 52).Value = 0.080211431263273
$ws.Cells.Item(190, 53).Value = 0.044686820246292
$ws.Cells.Item(190, 54).Value = 0.038138267885118
$ws.Cells.Item(190, 55).Value = 0.040404116881591
$ws.Cells.Item(190, 56).Value = 0.045041125130394
$ws.Cells.Item(190, 57).Value = 0.040687592972323

# Row 191
$ws.Cells.Item(191, 2).Value = 0.050008800847729
$ws.Cells.Item(191, 3).Value = 0.027570444089411
$ws.Cells.Item(191, 4).Value = 0.027122760834122
$ws.Cells.Item(191, 6).Value = 0.039554258653601
$ws.Cells.Item(191, 7).Value = 0.037106964252548
$ws.Cells.Item(191, 8).Value = 0.03951930651182
$ws.Cells.Item(191, 9).Value = 0.033985620573749
$ws.Cells.Item(191, 10).Value = 0.044871745947817
$ws.Cells.Item(191, 11).Value = 0.030801086197809
$ws.Cells.Item(191, 12).Value = 0.026549777382249
$ws.Cells.Item(191, 13).Value = 0.021574583511303
$ws.Cells.Item(191, 14).Value = 0.055084823300144
$ws.Cells.Item(191, 15).Value = 0.02262526695456
$ws.Cells.Item(191, 16).Value = 0.023248772209771
$ws.Cells.Item(191, 17).Value = 0.029319463764931
$ws.Cells.Item(191, 18).Value = 0.025923660164213
$ws.Cells.Item(191, 19).Value = 0.024186619928456
$ws.Cells.Item(191, 20).Value = 0.023513561673158
$ws.Cells.Item(191, 21).Value = 0.024284419172813
$ws.Cells.Item(191, 22).Value = 0.028879147381956
$ws.Cells.Item(191, 23).Value = 0.036055083073461
$ws.Cells.Item(191, 24).Value = 0.028620699371607
$ws.Cells.Item(191, 25).Value = 0.033172447482928
$ws.Cells.Item(191, 26).Value = 0.033440158388509
$ws.Cells.Item(191, 27).Value = 0.029424904410711
$ws.Cells.Item(191, 28).Value = 0.028165682699723
$ws.Cells.Item(191, 29).Value = 0.070679715157107
$ws.Cells.Item(191, 30).Value = 0.027332517839947
$ws.Cells.Item(191, 31).Value = 0.039598481318641
$ws.Cells.Item(191, 32).Value = 0.023854510003799
$ws.Cells.Item(191, 33).Value = 0.027913918175781
$ws.Cells.Item(191, 34).Value = 0.027057301540282
$ws.Cells.Item(191, 35).Value = 0.027154634713732
$ws.Cells.Item(191, 36).Value = 0.027215201823004
$ws.Cells.Item(191, 37).Value = 0.035388970697199
$ws.Cells.Item(191, 38).Value = 0.03220397694606
$ws.Cells.Item(191, 39).Value = 0.032148266777358
$ws.Cells.Item(191, 40).Value = 0.028619178290637
$ws.Cells.Item(191, 41).Value = 0.025267553490815
$ws.Cells.Item(191, 42).Value = 0.032596770557325
$ws.Cells.Item(191, 43).Value = 0.031133577107927
$ws.Cells.Item(191, 44).Value = 0.067009470831699
$ws.Cells.Item(191, 45).Value = 0.034581067914036
$ws.Cells.Item(191, 46).Value = 0.025855758816295
$ws.Cells.Item(191, 47).Value = 0.028578238898173
$ws.Cells.Item(191, 48).Value = 0.023737866634495
$ws.Cells.Item(191, 49).Value = 0.023711087084311
$ws.Cells.Item(191, 50).Value = 0.030261320577445
$ws.Cells.Item(191, 51).Value = 0.026531101606938
$ws.Cells.Item(191, 52).Value = 0.058976250307811
$ws.Cells.Item(191, 53).Value = 0.036136340651495
$ws.Cells.Item(191, 54).Value = 0.029782379248452
$ws.Cells.Item(191, 55).Value = 0.026604890578279
$ws.Cells.Item(191, 56).Value = 0.030928262368903
$ws.Cells.Item(191, 57).Value = 0.032280848460166

# Row 192
$ws.Cells.Item(192, 2).Value = 0.041858882248227
$ws.Cells.Item(192, 3).Value = 0.024943130335437
$ws.Cells.Item(192, 4).Value = 0.025418668566278
$ws.Cells.Item(192, 6).Value = 0.028893967083693
$ws.Cells.Item(192, 7).Value = 0.032308915086551
$ws.Cells.Item(192, 8).Value = 0.031065501848194
$ws.Cells.Item(192, 9).Value = 0.028228845941633
$ws.Cells.Item(192, 10).Value = 0.039087079242604
$ws.Cells.Item(192, 11).Value = 0.026368402163893
$ws.Cells.Item(192, 12).Value = 0.022726264269128
$ws.Cells.Item(192, 13).Value = 0.020978955022148
$ws.Cells.Item(192, 14).Value = 0.058794399957297
$ws.Cells.Item(192, 15).Value = 0.027162568076312
$ws.Cells.Item(192, 16).Value = 0.02843518174336
$ws.Cells.Item(192, 17).Value = 0.030683357068886
$ws.Cells.Item(192, 18).Value = 0.028589300177576
$ws.Cells.Item(192, 19).Value = 0.027314579664315
$ws.Cells.Item(192, 20).Value = 0.028815000614117
$ws.Cells.Item(192, 21).Value = 0.028471882771874
$ws.Cells.Item(192, 22).Value = 0.032743893898438
$ws.Cells.Item(192, 23).Value = 0.038732518019982
$ws.Cells.Item(192, 24).Value = 0.031685098971385
$ws.Cells.Item(192, 25).Value = 0.037737913362481
$ws.Cells.Item(192, 26).Value = 0.031185482184134
$ws.Cells.Item(192, 27).Value = 0.025351669328223
$ws.Cells.Item(192, 28).Value = 0.023745446004216
$ws.Cells.Item(192, 29).Value = 0.069819345214322
$ws.Cells.Item(192, 30).Value = 0.023120282704615
$ws.Cells.Item(192, 31).Value = 0.034944982478798
$ws.Cells.Item(192, 32).Value = 0.02670201280269
$ws.Cells.Item(192, 33).Value = 0.030904615101585
$ws.Cells.Item(192, 34).Value = 0.022322900482071
$ws.Cells.Item(192, 35).Value = 0.025074056425676
$ws.Cells.Item(192, 36).Value = 0.025728030842301
$ws.Cells.Item(192, 37).Value = 0.034815578943362
$ws.Cells.Item(192, 38).Value = 0.027752174291188
$ws.Cells.Item(192, 39).Value = 0.036529530598401
$ws.Cells.Item(192, 40).Value = 0.033251756796682
$ws.Cells.Item(192, 41).Value = 0.026701887184068
$ws.Cells.Item(192, 42).Value = 0.034290243220188
$ws.Cells.Item(192, 43).Value = 0.033710473508115
$ws.Cells.Item(192, 44).Value = 0.070757005470694
$ws.Cells.Item(192, 45).Value = 0.035794342229237
$ws.Cells.Item(192, 46).Value = 0.025267530748723
$ws.Cells.Item(192, 47).Value = 0.028986277280724
$ws.Cells.Item(192, 48).Value = 0.022798976519621
$ws.Cells.Item(192, 49).Value = 0.023139560346661
$ws.Cells.Item(192, 50).Value = 0.028998996096921
$ws.Cells.Item(192, 51).Value = 0.026593690881395
$ws.Cells.Item(192, 52).Value = 0.061461982636767
$ws.Cells.Item(192, 53).Value = 0.035897800142609
$ws.Cells.Item(192, 54).Value = 0.028968921293238
$ws.Cells.Item(192, 55).Value = 0.025344254831201
$ws.Cells.Item(192, 56).Value = 0.030071609036338
$ws.Cells.Item(192, 57).Value = 0.028516029311783

# Row 193
$ws.Cells.Item(193, 2).Value = 0.046026090423966
$ws.Cells.Item(193, 3).Value = 0.052918434301101
$ws.Cells.Item(193, 4).Value = 0.048665653599385
$ws.Cells.Item(193, 6).Value = 0.042049218024318
$ws.Cells.Item(193, 7).Value = 0.057601503803822
$ws.Cells.Item(193, 8).Value = 0.057970090860322
$ws.Cells.Item(193, 9).Value = 0.052501903248581
$ws.Cells.Item(193, 10).Value = 0.056572930298422
$ws.Cells.Item(193, 11).Value = 0.047155315545943
$ws.Cells.Item(193, 12).Value = 0.052537603497921
$ws.Cells.Item(193, 13).Value = 0.054216661876463
$ws.Cells.Item(193, 14).Value = 0.057088888859292
$ws.Cells.Item(193, 15).Value = 0.030437349843032
$ws.Cells.Item(193, 16).Value = 0.042919191157467
$ws.Cells.Item(193, 17).Value = 0.049294669596677
$ws.Cells.Item(193, 18).Value = 0.041788038232951
$ws.Cells.Item(193, 19).Value = 0.054070959840999
$ws.Cells.Item(193, 20).Value = 0.04796496136894
$ws.Cells.Item(193, 21).Value = 0.050266798195428
$ws.Cells.Item(193, 22).Value = 0.061393050688178
$ws.Cells.Item(193, 23).Value = 0.04923868601934
$ws.Cells.Item(193, 24).Value = 0.046749801948849
$ws.Cells.Item(193, 25).Value = 0.049684247867941
$ws.Cells.Item(193, 26).Value = 0.047246162348822
$ws.Cells.Item(193, 27).Value = 0.041583760255156
$ws.Cells.Item(193, 28).Value = 0.047797746021041
$ws.Cells.Item(193, 29).Value = 0.084799322019094
$ws.Cells.Item(193, 30).Value = 0.061060656702218
$ws.Cells.Item(193, 31).Value = 0.050574382264206
$ws.Cells.Item(193, 32).Value = 0.048757655011974
$ws.Cells.Item(193, 33).Value = 0.042470796960375
$ws.Cells.Item(193, 34).Value = 0.048998107317967
$ws.Cells.Item(193, 35).Value = 0.044129888656513
$ws.Cells.Item(193, 36).Value = 0.042124802938824
$ws.Cells.Item(193, 37).Value = 0.046142685530124
$ws.Cells.Item(193, 38).Value = 0.043724496674095
$ws.Cells.Item(193, 39).Value = 0.044130521094427
$ws.Cells.Item(193, 40).Value = 0.045877192368978
$ws.Cells.Item(193, 41).Value = 0.049147513985993
$ws.Cells.Item(193, 42).Value = 0.044226272142241
$ws.Cells.Item(193, 43).Value = 0.046766278906858
$ws.Cells.Item(193, 44).Value = 0.075668353268875
$ws.Cells.Item(193, 45).Value = 0.040208602977681
$ws.Cells.Item(193, 46).Value = 0.052756975802968
$ws.Cells.Item(193, 47).Value = 0.050516921132961
$ws.Cells.Item(193, 48).Value = 0.057726146634041
$ws.Cells.Item(193, 49).Value = 0.047750651658817
$ws.Cells.Item(193, 50).Value = 0.057106921774819
$ws.Cells.Item(193, 51).Value = 0.047774528803678
$ws.Cells.Item(193, 52).Value = 0.071064971630743
$ws.Cells.Item(193, 53).Value = 0.051729254728598
$ws.Cells.Item(193, 54).Value = 0.043803726043286
$ws.Cells.Item(193, 55).Value = 0.045604905986004
$ws.Cells.Item(193, 56).Value = 0.050528584427298
$ws.Cells.Item(193, 57).Value = 0.04709803828968

# Row 194
$ws.Cells.Item(194, 1).Value = "11 08 2020"
$ws.Cells.Item(194, 2).Value = 0.044118252638562
$ws.Cells.Item(194, 3).Value = 0.054033371340887
$ws.Cells.Item(194, 4).Value = 0.049744691606255
$ws.Cells.Item(194, 6).Value = 0.039724046796903
$ws.Cells.Item(194, 7).Value = 0.055069727957749
$ws.Cells.Item(194, 8).Value = 0.059806629495529
$ws.Cells.Item(194, 9).Value = 0.0498136456508
$ws.Cells.Item(194, 10).Value = 0.047877928037791
$ws.Cells.Item(194, 11).Value = 0.047246762628939
$ws.Cells.Item(194, 12).Value = 0.049326191772621
$ws.Cells.Item(194, 13).Value = 0.052845157262935
$ws.Cells.Item(194, 14).Value = 0.053758491328396
$ws.Cells.Item(194, 15).Value = 0.027599526246962
$ws.Cells.Item(194, 16).Value = 0.042808306786451
$ws.Cells.Item(194, 17).Value = 0.048714121240212
$ws.Cells.Item(194, 18).Value = 0.041579658624977
$ws.Cells.Item(194, 19).Value = 0.05368288136542
$ws.Cells.Item(194, 20).Value = 0.047990974298318
$ws.Cells.Item(194, 21).Value = 0.050121571272439
$ws.Cells.Item(194, 22).Value = 0.062020286741533
$ws.Cells.Item(194, 23).Value = 0.048708481728583
$ws.Cells.Item(194, 24).Value = 0.046091649209684
$ws.Cells.Item(194, 25).Value = 0.04908488442969
$ws.Cells.Item(194, 26).Value = 0.047615784766575
$ws.Cells.Item(194, 27).Value = 0.042221180892404
$ws.Cells.Item(194, 28).Value = 0.049903936170304
$ws.Cells.Item(194, 29).Value = 0.06756455861343399
$ws.Cells.Item(194, 30).Value = 0.064376748076772
$ws.Cells.Item(194, 31).Value = 0.053971932195427
$ws.Cells.Item(194, 32).Value = 0.048415953711232
$ws.Cells.Item(194, 33).Value = 0.042579852090841
$ws.Cells.Item(194, 34).Value = 0.051734084008053
$ws.Cells.Item(194, 35).Value = 0.046726427699178
$ws.Cells.Item(194, 36).Value = 0.043940341495793
$ws.Cells.Item(194, 37).Value = 0.047744744471779
$ws.Cells.Item(194, 38).Value = 0.045742128540185
$ws.Cells.Item(194, 39).Value = 0.044687595984735
$ws.Cells.Item(194, 40).Value = 0.044759871571587
$ws.Cells.Item(194, 41).Value = 0.05015239809719
$ws.Cells.Item(194, 42).Value = 0.044127364913003
$ws.Cells.Item(194, 43).Value = 0.049133615336722
$ws.Cells.Item(194, 44).Value = 0.076293671624287
$ws.Cells.Item(194, 45).Value = 0.05325480958244
$ws.Cells.Item(194, 46).Value = 0.056825537801944
$ws.Cells.Item(194, 47).Value = 0.055287311367109
$ws.Cells.Item(194, 48).Value = 0.061968473349256
$ws.Cells.Item(194, 49).Value = 0.05392668128226
$ws.Cells.Item(194, 50).Value = 0.063826011111523
$ws.Cells.Item(194, 51).Value = 0.05366975045501
$ws.Cells.Item(194, 52).Value = 0.069287972896525
$ws.Cells.Item(194, 53).Value = 0.05655784024162
$ws.Cells.Item(194, 54).Value = 0.048153483545589
$ws.Cells.Item(194, 55).Value = 0.049604241617066
$ws.Cells.Item(194, 56).Value = 0.054514876456383
$ws.Cells.Item(194, 57).Value = 0.048635028357488

# Row 195
$ws.Cells.Item(195, 1).Value = "12 08 2020"
$ws.Cells.Item(195, 2).Value = 0.051835399302178
$ws.Cells.Item(195, 3).Value = 0.054676532820954
$ws.Cells.Item(195, 4).Value = 0.052631227119385
$ws.Cells.Item(195, 6).Value = 0.047919724265789
$ws.Cells.Item(195, 7).Value = 0.046744429885148
$ws.Cells.Item(195, 8).Value = 0.054573775563135
$ws.Cells.Item(195, 9).Value = 0.047738010729529
$ws.Cells.Item(195, 10).Value = 0.045435787487791
$ws.Cells.Item(195, 11).Value = 0.04562891135575
$ws.Cells.Item(195, 12).Value = 0.049440157113606
$ws.Cells.Item(195, 13).Value = 0.052556729410555
$ws.Cells.Item(195, 14).Value = 0.059418342479345
$ws.Cells.Item(195, 15).Value = 0.031956258849025
$ws.Cells.Item(195, 16).Value = 0.046294215766254
$ws.Cells.Item(195, 17).Value = 0.050749651555072
$ws.Cells.Item(195, 18).Value = 0.04240100999012
$ws.Cells.Item(195, 19).Value = 0.055769746525046
$ws.Cells.Item(195, 20).Value = 0.049172059753568
$ws.Cells.Item(195, 21).Value = 0.050883229563625
$ws.Cells.Item(195, 22).Value = 0.061223516042335
$ws.Cells.Item(195, 23).Value = 0.048608644386112
$ws.Cells.Item(195, 24).Value = 0.045893481769109
$ws.Cells.Item(195, 25).Value = 0.050223851035317
$ws.Cells.Item(195, 26).Value = 0.046875794628247
$ws.Cells.Item(195, 27).Value = 0.04360476057426
$ws.Cells.Item(195, 28).Value = 0.050001102320727
$ws.Cells.Item(195, 29).Value = 0.073068237287524
$ws.Cells.Item(195, 30).Value = 0.061496368484489
$ws.Cells.Item(195, 31).Value = 0.053032797245249
$ws.Cells.Item(195, 32).Value = 0.050800613276599
$ws.Cells.Item(195, 33).Value = 0.04719659443808
$ws.Cells.Item(195, 34).Value = 0.057500292436793
$ws.Cells.Item(195, 35).Value = 0.048520529981828
$ws.Cells.Item(195, 36).Value = 0.046430310682728
$ws.Cells.Item(195, 37).Value = 0.04974853892368
$ws.Cells.Item(195, 38).Value = 0.047914605686715
$ws.Cells.Item(195, 39).Value = 0.04887763383805
$ws.Cells.Item(195, 40).Value = 0.047378808740902
$ws.Cells.Item(195, 41).Value = 0.051334125746039
$ws.Cells.Item(195, 42).Value = 0.044812672799498
$ws.Cells.Item(195, 43).Value = 0.046775538145325
$ws.Cells.Item(195, 44).Value = 0.079470295328549
$ws.Cells.Item(195, 45).Value = 0.047795887614522
$ws.Cells.Item(195, 46).Value = 0.0513426811407
$ws.Cells.Item(195, 47).Value = 0.048230227395261
$ws.Cells.Item(195, 48).Value = 0.055959009775526
$ws.Cells.Item(195, 49).Value = 0.047262728432198
$ws.Cells.Item(195, 50).Value = 0.057070857206656
$ws.Cells.Item(195, 51).Value = 0.04651625296955
$ws.Cells.Item(195, 52).Value = 0.085232114105792
$ws.Cells.Item(195, 53).Value = 0.050280341499272
$ws.Cells.Item(195, 54).Value = 0.044325902421577
$ws.Cells.Item(195, 55).Value = 0.047945950481005
$ws.Cells.Item(195, 56).Value = 0.050635588292848
$ws.Cells.Item(195, 57).Value = 0.050467609879462

# Row 196
$ws.Cells.Item(196, 1).Value = "13 08 2020"

# Row 197
$ws.Cells.Item(197, 1).Value = "14 08 2020"

# Row 198
$ws.Cells.Item(198, 1).Value = "15 08 2020"

# Row 199
$ws.Cells.Item(199, 1).Value = "16 08 2020"

# Row 200
$ws.Cells.Item(200, 1).Value = "17 08 2020"
